$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "attributes" sheet: insert a new "label" column (B) right after "name".
#    By default label = name, except for the "Regexp" attribute of
#    SAM_TAG_Types, which is renamed to "Regexp_" while its label stays
#    "Regexp" (fixes an import collision, see commit message).
# ---------------------------------------------------------------------------
$attrs = $wb.Worksheets.Item("attributes")

# Insert new column before current column B ("entity").
$attrs.Range("B1").EntireColumn.Insert()

# New header for inserted column.
$attrs.Range("B1").Value = "label"

# Populate label column with the attribute name for every data row (2-20).
for ($r = 2; $r -le 20; $r++) {
    $nameCell = $attrs.Cells.Item($r, 1)
    $attrs.Cells.Item($r, 2).Value = $nameCell.Value2
}

# Row 19 used to be the "Regexp" attribute of SAM_TAG_Types; rename it to
# "Regexp_" and keep "Regexp" as its label.
$attrs.Range("A19").Value = "Regexp_"
$attrs.Range("B19").Value = "Regexp"

$attrs.Range("A19").Select()

# ---------------------------------------------------------------------------
# 2. "SAM_TAG_Types" sheet: the data column header must follow the renamed
#    attribute.
# ---------------------------------------------------------------------------
$tagTypes = $wb.Worksheets.Item("SAM_TAG_Types")
$tagTypes.Range("B1").Value = "Regexp_"
$tagTypes.Range("D12").Select()
